# Updated cryptos list on Tue Jun 13 17:43:00 UTC 2023 with GitHub Actions
# This script rewrites the Coin / Link / Price / Volume(1h) table (rows 2-51)
# to reflect the latest scrape. Column A (rank index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "25.843.64", "  +0.19%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.736.88", "  +0.32%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.001", "  -0.03%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "236.79", "  +3.27%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  -0.03%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5119", "  -0.70%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2736", "  -0.23%  ")
    ,@(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "39.99", "  +1.90%  ")
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06105", "  -0.01%  ")
    ,@(11, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.740.02", "  +0.34%  ")
    ,@(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07171", "  +1.91%  ")
    ,@(13, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "14.91", "  +0.41%  ")
    ,@(14, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6358", "  +0.93%  ")
    ,@(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.589", "  +2.10%  ")
    ,@(16, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "77.16", "  +1.05%  ")
    ,@(17, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.01%  ")
    ,@(18, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  -0.01%  ")
    ,@(19, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "25.852.76", "  +0.17%  ")
    ,@(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.70", "  +2.73%  ")
    ,@(21, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000006744", "  +2.31%  ")
    ,@(22, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.959.13", "  +0.24%  ")
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.255", "  +2.56%  ")
    ,@(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.669", "  -0.22%  ")
    ,@(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.222", "  +1.60%  ")
    ,@(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "138.95", "  -0.56%  ")
    ,@(27, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.513", "  +0.62%  ")
    ,@(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "15.11", "  +0.82%  ")
    ,@(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.752", "  -0.13%  ")
    ,@(30, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "105.74", "  +3.94%  ")
    ,@(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.003", "  +9.30%  ")
    ,@(32, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08318", "  +0.61%  ")
    ,@(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.644", "  +5.14%  ")
    ,@(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04563", "  +2.38%  ")
    ,@(35, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.668", "  +2.14%  ")
    ,@(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9844", "  +1.68%  ")
    ,@(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6159", "  +1.04%  ")
    ,@(38, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.686", "  +1.30%  ")
    ,@(39, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01608", "  +2.55%  ")
    ,@(40, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.918", "  +0.61%  ")
    ,@(41, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.001", "  -0.07%  ")
    ,@(42, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "97.91", "  -1.34%  ")
    ,@(43, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3832", "  +1.05%  ")
    ,@(44, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.7356", "  +2.57%  ")
    ,@(45, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "4.949", "  -0.57%  ")
    ,@(46, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1122", "  -0.03%  ")
    ,@(47, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05258", "  -1.67%  ")
    ,@(48, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.152", "  +0.28%  ")
    ,@(49, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "54.77", "  +3.86%  ")
    ,@(50, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "30.49", "  +2.31%  ")
    ,@(51, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.535", "  +0.13%  ")
)

foreach ($item in $data) {
    $row   = $item[0]
    $coin  = $item[1]
    $link  = $item[2]
    $price = $item[3]
    $vol   = $item[4]

    # Coin name and link are plain (non-numeric-looking) text - safe to set directly.
    $ws.Cells.Item($row, 2).Value = $coin
    $ws.Cells.Item($row, 3).Value = $link

    # Price values often look like plain numbers (e.g. "1.001", "236.79") which
    # Excel would otherwise silently convert to a numeric type. Force the cell
    # to text first, assign the literal string, then clear the number-format
    # override again so the cell keeps the workbook's original (default) style.
    $priceCell = $ws.Cells.Item($row, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.ClearFormats()

    # Volume/percentage text (e.g. "  +0.19%  ") is never numeric-looking, but
    # use the same safe approach for consistency and to preserve the
    # leading/trailing spaces exactly.
    $volCell = $ws.Cells.Item($row, 5)
    $volCell.NumberFormat = "@"
    $volCell.Value = $vol
    $volCell.ClearFormats()
}

Write-Host "Updated $($data.Count) crypto rows"
